$wb = $excel.ActiveWorkbook

# Update the status text for the "a5602cf4-..." report row from
# "Ready for handoff" to "Handback transform failed" everywhere it appears
# (Overview sheet columns B3/C3, and the "Status" column C3 on the zh-cn
# and de-de report sheets all share this string).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Record the "Error Detail" message produced by the failed handback
# transform for each locale's report row.
$wsZhCn.Range("K3").Value = "Handback file name: jxadawnf.521 is different with handoff file name: a5602cf4-2330-46dc-8709-2b8558240a21.9a1b0216089abdf0d8056bfb8fdfd5c4c4800a17.zh-cn."

$wsDeDe.Range("K3").Value = "Handback file name: jxadawnf.521 is different with handoff file name: a5602cf4-2330-46dc-8709-2b8558240a21.9a1b0216089abdf0d8056bfb8fdfd5c4c4800a17.de-de."
